# Updates cryptos list with latest price/volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values. These look numeric (e.g. "511.03", "1.00") so a
# literal Value assignment would make Excel auto-coerce them to a Number.
# Prefixing with an apostrophe forces literal text (as Excel's UI does),
# then resetting Style back to "Normal" clears the quotePrefix formatting
# Excel applies, leaving the cell's style untouched.
$priceUpdates = @{
    'D2' = '56.740.07';
    'D3' = '3.028.47';
    'D5' = '511.03';
    'D6' = '140.21';
    'D9' = '7.15';
    'D11' = '0.370';
    'D12' = '3.540.46';
    'D14' = '25.30';
    'D15' = '0.0000163';
    'D16' = '56.720.07';
    'D17' = '3.026.65';
    'D18' = '5.93';
    'D19' = '13.17';
    'D20' = '8.04';
    'D21' = '333.99';
    'D22' = '1.00';
    'D24' = '64.78';
    'D25' = '3.149.86';
    'D27' = '1.00';
    'D29' = '6.41';
    'D33' = '20.44';
    'D34' = '153.03';
    'D35' = '4.47';
    'D36' = '26.92';
    'D40' = '3.065.55';
    'D41' = '36.57';
    'D42' = '1.00';
    'D44' = '0.658';
    'D45' = '2.205.26';
    'D48' = '0.931';
    'D50' = '19.75';
    'D51' = '0.0856'
}

foreach ($ref in $priceUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.Value = "'" + $priceUpdates[$ref]
    $cell.Style = "Normal"
}

# Column E ("Volume(1h)") values already contain surrounding spaces, so they
# are never misread as numbers and can be assigned directly.
$volumeUpdates = @{
    'E2' = '  +0.20%  ';
    'E3' = '  +2.38%  ';
    'E4' = '  -0.12%  ';
    'E5' = '  +2.94%  ';
    'E6' = '  +4.29%  ';
    'E7' = '  -0.05%  ';
    'E8' = '  +1.43%  ';
    'E9' = '  -0.32%  ';
    'E10' = '  +2.30%  ';
    'E11' = '  +5.00%  ';
    'E12' = '  +2.22%  ';
    'E13' = '  +0.57%  ';
    'E14' = '  -2.46%  ';
    'E15' = '  +3.86%  ';
    'E16' = '  +0.02%  ';
    'E17' = '  +2.32%  ';
    'E18' = '  -0.76%  ';
    'E19' = '  +5.76%  ';
    'E20' = '  +3.83%  ';
    'E21' = '  +5.44%  ';
    'E22' = '  +0.17%  ';
    'E23' = '  +3.26%  ';
    'E24' = '  +3.48%  ';
    'E25' = '  +2.28%  ';
    'E26' = '  +2.75%  ';
    'E27' = '  +0.60%  ';
    'E28' = '  +8.70%  ';
    'E29' = '  -1.03%  ';
    'E30' = '  -3.06%  ';
    'E32' = '  +3.59%  ';
    'E33' = '  +2.76%  ';
    'E34' = '  +0.07%  ';
    'E35' = '  +0.06%  ';
    'E36' = '  +12.85%  ';
    'E37' = '  +2.50%  ';
    'E38' = '  +1.67%  ';
    'E39' = '  +1.43%  ';
    'E40' = '  +2.55%  ';
    'E41' = '  -2.02%  ';
    'E42' = '  -0.19%  ';
    'E43' = '  +3.30%  ';
    'E44' = '  +2.78%  ';
    'E45' = '  +2.40%  ';
    'E46' = '  +0.12%  ';
    'E47' = '  +4.73%  ';
    'E48' = '  +0.99%  ';
    'E49' = '  +0.14%  ';
    'E50' = '  +4.21%  ';
    'E51' = '  -0.02%  '
}

foreach ($ref in $volumeUpdates.Keys) {
    $ws.Range($ref).Value = $volumeUpdates[$ref]
}
